# This workbook holds a weekly price-report table (Cebollín / Terminal La
# Palmera de La Serena). A new week's record is inserted at the top of the
# data block (row 78), pushing every existing record in that block down by
# one row; the oldest record in the block ends up duplicated one row lower
# (i.e. the table simply grows by one row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 78; Excel shifts rows 78:170 down to 79:171
# (and the entire row's formatting, incl. the date-formatted column D style).
$ws.Rows.Item(78).Insert()

# Populate the freshly inserted row 78 with the new weekly record.
$ws.Range("A78").Value = 8
$ws.Range("B78").Value = "Terminal La Palmera de La Serena"
$ws.Range("C78").Value = "Coquimbo"
$ws.Range("D78").Value = 44638
$ws.Range("E78").Value = 4
$ws.Range("F78").Value = 100112037
$ws.Range("G78").Value = "Cebollín"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 1800
$ws.Range("K78").Value = 1100
$ws.Range("L78").Value = 1200
$ws.Range("M78").Value = 1150
$ws.Range("N78").Value = "$/paquete 6 unidades"
$ws.Range("O78").Value = "Provincia del Elquí"
$ws.Range("P78").Value = 192
$ws.Range("Q78").Value = 6
$ws.Range("R78").Value = "Hortaliza"
